# "Restore front page elements"
#
# The author restored the front-page view of the glossary sheet:
#   - widened column B (the "definition" column) back to its full width
#   - restored the active window position/size
#   - restored the selected cell to C11 (bottom of the list) instead of I12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore column B's wide custom width (~98.16 characters, holding long
# glossary definitions). ColumnWidth is specified in character units; Excel
# stores it (in the OOXML <col width>) offset by the default column padding.
$ws.Columns.Item(2).ColumnWidth = 97.3333333333

# Restore the document window's on-screen position/size.
$win = $excel.ActiveWindow
$win.Left = 10440
$win.Top = 460
$win.Width = 23820
$win.Height = 16140

# Restore the saved selection/active cell to C11.
$ws.Range("C11").Select() | Out-Null
